$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.9999674344715328,
    0.9989400190852257,
    0.9999619811694375,
    0.9999463330267169,
    0.9999562910719529,
    0.00003039849136261971,
    0.0009894456561571146,
    0.00003759280563521777,
    0.00005489286478983781,
    0.0000462428352125278,
    0.000349036811621329,
    0.005513482689065026,
    1.000060120975632,
    0.005748203061621094,
    94.80223515402628,
    139.9006406741497
)

for ($row = 2; $row -le 26; $row++) {
    for ($i = 0; $i -lt $newValues.Count; $i++) {
        $col = $i + 2  # Column B = 2
        $ws.Cells.Item($row, $col).Value = $newValues[$i]
    }
}
